# Weekly price-sheet update: a new Cilantro price record for
# "Feria Lagunitas de Puerto Montt" is inserted as row 92, pushing the
# existing rows 92-176 down to rows 93-177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92 (shifts rows 92..176 down to 93..177).
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record's data.
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 44484
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100112040
$ws.Cells.Item(92, 7).Value = "Cilantro"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 300
$ws.Cells.Item(92, 11).Value = 11000
$ws.Cells.Item(92, 12).Value = 11500
$ws.Cells.Item(92, 13).Value = 11250
$ws.Cells.Item(92, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 312
$ws.Cells.Item(92, 17).Value = 36
$ws.Cells.Item(92, 18).Value = "Hortaliza"
